$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet1 (BD_Times): append rows 492-505 ---
$ws1.Cells.Item(492,1).Value = "Chapecoense"
$ws1.Cells.Item(492,2).Value = 1
$ws1.Cells.Item(492,3).Value = 0
$ws1.Cells.Item(492,4).Value = 0
$ws1.Cells.Item(492,5).Value = 0
$ws1.Cells.Item(492,6).Value = 0
$ws1.Cells.Item(492,7).Value = 0
$ws1.Cells.Item(492,8).Value = 10
$ws1.Cells.Item(492,9).Value = 5

$ws1.Cells.Item(493,1).Value = "Avai"
$ws1.Cells.Item(493,2).Value = 0
$ws1.Cells.Item(493,3).Value = 0
$ws1.Cells.Item(493,4).Value = 0
$ws1.Cells.Item(493,5).Value = 0
$ws1.Cells.Item(493,6).Value = 0
$ws1.Cells.Item(493,7).Value = 0
$ws1.Cells.Item(493,8).Value = 5
$ws1.Cells.Item(493,9).Value = 10

$ws1.Cells.Item(494,1).Value = "Atletico-GO"
$ws1.Cells.Item(494,2).Value = 1
$ws1.Cells.Item(494,3).Value = 0
$ws1.Cells.Item(494,4).Value = 0
$ws1.Cells.Item(494,5).Value = 0
$ws1.Cells.Item(494,6).Value = 0
$ws1.Cells.Item(494,7).Value = 0
$ws1.Cells.Item(494,8).Value = 4
$ws1.Cells.Item(494,9).Value = 3

$ws1.Cells.Item(495,1).Value = "Vitoria"
$ws1.Cells.Item(495,2).Value = 0
$ws1.Cells.Item(495,3).Value = 0
$ws1.Cells.Item(495,4).Value = 0
$ws1.Cells.Item(495,5).Value = 0
$ws1.Cells.Item(495,6).Value = 0
$ws1.Cells.Item(495,7).Value = 0
$ws1.Cells.Item(495,8).Value = 3
$ws1.Cells.Item(495,9).Value = 4

$ws1.Cells.Item(496,1).Value = "Botafogo"
$ws1.Cells.Item(496,2).Value = 1
$ws1.Cells.Item(496,3).Value = 0
$ws1.Cells.Item(496,4).Value = 0
$ws1.Cells.Item(496,5).Value = 0
$ws1.Cells.Item(496,6).Value = 0
$ws1.Cells.Item(496,7).Value = 0
$ws1.Cells.Item(496,8).Value = 6
$ws1.Cells.Item(496,9).Value = 3

$ws1.Cells.Item(497,1).Value = "ABC"
$ws1.Cells.Item(497,2).Value = 0
$ws1.Cells.Item(497,3).Value = 0
$ws1.Cells.Item(497,4).Value = 0
$ws1.Cells.Item(497,5).Value = 0
$ws1.Cells.Item(497,6).Value = 0
$ws1.Cells.Item(497,7).Value = 0
$ws1.Cells.Item(497,8).Value = 3
$ws1.Cells.Item(497,9).Value = 6

$ws1.Cells.Item(498,1).Value = "Novohorizontino"
$ws1.Cells.Item(498,2).Value = 1
$ws1.Cells.Item(498,3).Value = 1
$ws1.Cells.Item(498,4).Value = 0
$ws1.Cells.Item(498,5).Value = 0
$ws1.Cells.Item(498,6).Value = 2
$ws1.Cells.Item(498,7).Value = 0
$ws1.Cells.Item(498,8).Value = 4
$ws1.Cells.Item(498,9).Value = 4

$ws1.Cells.Item(499,1).Value = "Ituano"
$ws1.Cells.Item(499,2).Value = 0
$ws1.Cells.Item(499,3).Value = 0
$ws1.Cells.Item(499,4).Value = 1
$ws1.Cells.Item(499,5).Value = 0
$ws1.Cells.Item(499,6).Value = 0
$ws1.Cells.Item(499,7).Value = 2
$ws1.Cells.Item(499,8).Value = 4
$ws1.Cells.Item(499,9).Value = 4

$ws1.Cells.Item(500,1).Value = "Juventude"
$ws1.Cells.Item(500,2).Value = 1
$ws1.Cells.Item(500,3).Value = 0
$ws1.Cells.Item(500,4).Value = 0
$ws1.Cells.Item(500,5).Value = 0
$ws1.Cells.Item(500,6).Value = 0
$ws1.Cells.Item(500,7).Value = 0
$ws1.Cells.Item(500,8).Value = 8
$ws1.Cells.Item(500,9).Value = 4

$ws1.Cells.Item(501,1).Value = "Chapecoense"
$ws1.Cells.Item(501,2).Value = 0
$ws1.Cells.Item(501,3).Value = 0
$ws1.Cells.Item(501,4).Value = 0
$ws1.Cells.Item(501,5).Value = 0
$ws1.Cells.Item(501,6).Value = 0
$ws1.Cells.Item(501,7).Value = 0
$ws1.Cells.Item(501,8).Value = 4
$ws1.Cells.Item(501,9).Value = 8

$ws1.Cells.Item(502,1).Value = "Vitoria"
$ws1.Cells.Item(502,2).Value = 1
$ws1.Cells.Item(502,3).Value = 0
$ws1.Cells.Item(502,4).Value = 0
$ws1.Cells.Item(502,5).Value = 0
$ws1.Cells.Item(502,6).Value = 0
$ws1.Cells.Item(502,7).Value = 0
$ws1.Cells.Item(502,8).Value = 6
$ws1.Cells.Item(502,9).Value = 4

$ws1.Cells.Item(503,1).Value = "Mirassol"
$ws1.Cells.Item(503,2).Value = 0
$ws1.Cells.Item(503,3).Value = 0
$ws1.Cells.Item(503,4).Value = 0
$ws1.Cells.Item(503,5).Value = 0
$ws1.Cells.Item(503,6).Value = 0
$ws1.Cells.Item(503,7).Value = 0
$ws1.Cells.Item(503,8).Value = 4
$ws1.Cells.Item(503,9).Value = 6

$ws1.Cells.Item(504,1).Value = "Londrina"
$ws1.Cells.Item(504,2).Value = 1
$ws1.Cells.Item(504,3).Value = 1
$ws1.Cells.Item(504,4).Value = 1
$ws1.Cells.Item(504,5).Value = 1
$ws1.Cells.Item(504,6).Value = 1
$ws1.Cells.Item(504,7).Value = 1
$ws1.Cells.Item(504,8).Value = 14
$ws1.Cells.Item(504,9).Value = 1

$ws1.Cells.Item(505,1).Value = "Tombense"
$ws1.Cells.Item(505,2).Value = 0
$ws1.Cells.Item(505,3).Value = 1
$ws1.Cells.Item(505,4).Value = 1
$ws1.Cells.Item(505,5).Value = 1
$ws1.Cells.Item(505,6).Value = 1
$ws1.Cells.Item(505,7).Value = 1
$ws1.Cells.Item(505,8).Value = 1
$ws1.Cells.Item(505,9).Value = 14

# --- Sheet2 (BD_Jogo): append rows 247-253 ---
$ws2.Cells.Item(247,1).Value = 0
$ws2.Cells.Item(247,2).Value = 0
$ws2.Cells.Item(247,3).Value = 15
$ws2.Cells.Item(247,4).Value = "Chapecoense"
$ws2.Cells.Item(247,5).Value = "Avai"

$ws2.Cells.Item(248,1).Value = 0
$ws2.Cells.Item(248,2).Value = 0
$ws2.Cells.Item(248,3).Value = 7
$ws2.Cells.Item(248,4).Value = "Atletico-GO"
$ws2.Cells.Item(248,5).Value = "Vitoria"

$ws2.Cells.Item(249,1).Value = 0
$ws2.Cells.Item(249,2).Value = 0
$ws2.Cells.Item(249,3).Value = 9
$ws2.Cells.Item(249,4).Value = "Botafogo"
$ws2.Cells.Item(249,5).Value = "ABC"

$ws2.Cells.Item(250,1).Value = 0
$ws2.Cells.Item(250,2).Value = 2
$ws2.Cells.Item(250,3).Value = 8
$ws2.Cells.Item(250,4).Value = "Novohorizontino"
$ws2.Cells.Item(250,5).Value = "Ituano"

$ws2.Cells.Item(251,1).Value = 0
$ws2.Cells.Item(251,2).Value = 0
$ws2.Cells.Item(251,3).Value = 12
$ws2.Cells.Item(251,4).Value = "Juventude"
$ws2.Cells.Item(251,5).Value = "Chapecoense"

$ws2.Cells.Item(252,1).Value = 0
$ws2.Cells.Item(252,2).Value = 0
$ws2.Cells.Item(252,3).Value = 10
$ws2.Cells.Item(252,4).Value = "Vitoria"
$ws2.Cells.Item(252,5).Value = "Mirassol"

$ws2.Cells.Item(253,1).Value = 1
$ws2.Cells.Item(253,2).Value = 2
$ws2.Cells.Item(253,3).Value = 15
$ws2.Cells.Item(253,4).Value = "Londrina"
$ws2.Cells.Item(253,5).Value = "Tombense"

